$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell text assignments, in the exact order they were originally entered
# (this reproduces the shared-string table order / indices 449-474).
$entries = @(
  @{Row=171; Text="The testing is purely manual at this level of maturity. Teams rely on manual testers to find defects and bugs. Teams need to start writing some automated tests"},
  @{Row=172; Text="The cycle time of features and functionality is impacted due to testing late in the lifecycle. Functionality is stopped from being deployed due to defects/bugs being found at later stages of the lifecycle."},
  @{Row=173; Text="Non-functional testing is considered as an afterthought when the Software is about to be released. Teams do not factor non-functional testing in their designs and implementation. Teams need to start considering non-functional testing such as how the application will behave under load into consideration."},
  @{Row=174; Text="There is a separate testing team that designs test cases based on the requirements document without interacting with the teams. There should be closer collaboration with the development team."},
  @{Row=175; Text="The code base has very little or no unit/integration tests. Developers rely on testers to find bugs/defects in their implementation."},
  @{Row=176; Text="The Quality Assurance assets are non-reusable due to it being very manual and specific to each application. The test cases are driven by the requirements so no frameworks in place. The team should start investigating frameworks, which can be extended for other test cases so that any team member is able to contribute to testing."},
  @{Row=177; Text="At this maturity level, some of the regression tests are automated eventhough all the business critical cases do not have automated regression tests."},
  @{Row=178; Text="The regression tests are run on an ad hoc basis as opposed to on every critical change to the application. The regression tests should be run as part of the build and release towards production automatically to move from this level"},
  @{Row=179; Text="The test cases are designed by the team working on the application in collaboration with the stake holders. The team is aware of the business critical functionality and designs test cases to cater for these."},
  @{Row=180; Text="Regression test packs are fully automated for the application."},
  @{Row=181; Text="Regression tests are triggered by the build pipeline and run as part of the normal pipeline process. They are scheduled as daily run tasks due to the typical slowness in regression tests execution."},
  @{Row=183; Text="The team has implemented smoke tests and these are run as part of every code check-in. The smoke tests report are also displayed on the team tests metrics on the dashboard."},
  @{Row=184; Text="Performance tests are automated and the performance metrics are displayed on the dashboards. The build pipeline is failed if the performance metrics expected is not met. Due to the slow nature of performance tests, these are run as scheduled tasks on a nightly basis."},
  @{Row=185; Text="Non-functional requirements such as response times are defined and measured with test cases that are also displayed on the dashboards."},
  @{Row=186; Text="The team designs the test cases as part of the iteration and when building the functionality as opposed to as a pre-project process. Test cases are adapted to suit changes in the product vision and goal."},
  @{Row=188; Text="Infrastructure is automated and can be generated from scripts in the version control management system. The infrastructure automation is also tested using tools like ChefSpec."},
  @{Row=189; Text="Automation security testing is implemented and is part of the build pipeline using tools such as ZAP (Zed-Attack Proxy). The build is failed if security testing fails."},
  @{Row=190; Text="Test are automated and only the exploratory type tests are left unautomated."},
  @{Row=191; Text="Quality metrics are measured and tracked. Historical trends of the quality metrics are displayed on dashboards to the teams to drive the right culture of always improving quality. "},
  @{Row=187; Text="The team writes unit tests using established frameworks and methodologies. The unit tests coverage is not less than 50% and the build is failed if the unit tests coverage is lower than this. Developers write tests as they build functionality. Testing is not out-sourced to quality assurance engineers."},
  @{Row=192; Text="Any team member can execute the tests via script runs or on the build pipelines."},
  @{Row=193; Text="There are documented processes in place to understand the root cause of test failures and how to respond to them especially for the non-functional requirements testing."},
  @{Row=194; Text="Tests artefacts have same significance and importance as the application code and are continually refactored and maintained."},
  @{Row=182; Text="At least 50% of the code base is test covered. The test coverage metrics are captured on the team pipeline dashboards and the build is failed if the coverage threshold is less than 50%. A failing tests is always fixed before extra functionality is built."},
  @{Row=195; Text="Performance metrics are clearly defined and available to the team members early in the lifecycle. These are displayed on the team walls for all team members to be aware of and reviewed every iteration."},
  @{Row=196; Text="Tests drive the release readiness while business drives the release decisions. A release candidate is only generated and available in the pipeline if all the tests have passed."}
)

foreach ($e in $entries) {
  $ws.Range("E" + $e.Row).Value = $e.Text
}

# Row heights (auto-fit-equivalent explicit heights) for the newly wrapped rows.
$heights = @(
  @{Row=171; Height=43.2},
  @{Row=172; Height=43.2},
  @{Row=173; Height=72},
  @{Row=174; Height=43.2},
  @{Row=175; Height=28.8},
  @{Row=176; Height=72},
  @{Row=177; Height=43.2},
  @{Row=178; Height=57.6},
  @{Row=179; Height=43.2},
  @{Row=181; Height=43.2},
  @{Row=182; Height=57.6},
  @{Row=183; Height=43.2},
  @{Row=184; Height=57.6},
  @{Row=185; Height=28.8},
  @{Row=186; Height=43.2},
  @{Row=187; Height=72},
  @{Row=188; Height=43.2},
  @{Row=189; Height=43.2},
  @{Row=190; Height=28.8},
  @{Row=191; Height=43.2},
  @{Row=192; Height=28.8},
  @{Row=193; Height=43.2},
  @{Row=194; Height=28.8},
  @{Row=195; Height=43.2},
  @{Row=196; Height=43.2}
)

foreach ($h in $heights) {
  $ws.Rows.Item($h.Row).RowHeight = $h.Height
}

# Restore the final selection to match where editing left off.
[void]$ws.Range("E197").Select()
